$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value while forcing text storage (so purely-numeric-looking
# strings like confidence scores "0.76" or bbox coordinate lists "962,713,1006,765"
# are not auto-converted into numbers by Excel's type inference).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value2 = $value
}

# Rows 7-11: values rotate "up" by one row (old row8 -> row7, old row9 -> row8,
# old row10 -> row9, old row11 -> row10, old row7 -> row11).
# Apply directly using the known old/new literal values from the diff.

# Row 7
$ws.Range("A7").Value2 = "283b6eda-9c83-4cdd-9524-c7c394f2dc89"
$ws.Range("D7").Value2 = "image_20250728214139_ppp0.jpg"
$ws.Range("E7").Value2 = "PLACA_20250717165933"
$ws.Range("F7").Value2 = "Beja"
$ws.Range("G7").Value2 = 38.02035
$ws.Range("H7").Value2 = -7.94715
Set-TextValue $ws.Range("I7") "962,713,1006,765"
Set-TextValue $ws.Range("J7") "0.76"

# Row 8
$ws.Range("A8").Value2 = "a19b65d1-6f97-4841-9e1c-7446a9be92b6"
Set-TextValue $ws.Range("I8") "967,614,1002,659"
Set-TextValue $ws.Range("J8") "0.73"

# Row 9
$ws.Range("A9").Value2 = "4be1b1cf-d480-453e-b5fb-d4ecd6764c4d"
Set-TextValue $ws.Range("I9") "702,633,740,690"
Set-TextValue $ws.Range("J9") "0.72"

# Row 10
$ws.Range("A10").Value2 = "dfd476d4-7689-4671-a076-78fe3ce806bb"
Set-TextValue $ws.Range("I10") "1254,850,1294,895"
Set-TextValue $ws.Range("J10") "0.67"

# Row 11
$ws.Range("A11").Value2 = "2117575c-4ae1-458c-b88a-fc40f40debdb"
$ws.Range("D11").Value2 = "image_20250727074723_ppp0.jpg"
$ws.Range("E11").Value2 = "PLACA_20250723145134"
$ws.Range("F11").Value2 = "Moura"
$ws.Range("G11").Value2 = 38.06587
$ws.Range("H11").Value2 = -7.221796
Set-TextValue $ws.Range("I11") "1490,161,1563,258"
Set-TextValue $ws.Range("J11") "0.62"

# Row 16: image filename + bounding box refinement
$ws.Range("D16").Value2 = "image_20250807111026_ppp0.jpg"
Set-TextValue $ws.Range("I16") "641,529,688,576"

# Row 17: image filename + bounding box + confidence refinement
$ws.Range("D17").Value2 = "image_20250807111026_ppp0.jpg"
Set-TextValue $ws.Range("I17") "793,481,831,526"
Set-TextValue $ws.Range("J17") "0.70"

# Row 18: image filename + bounding box + confidence refinement
$ws.Range("D18").Value2 = "image_20250808100711_ppp0.jpg"
Set-TextValue $ws.Range("I18") "1182,409,1232,451"
Set-TextValue $ws.Range("J18") "0.75"

$wb.Save()
